# Update the "想去人数" (interested-people count) column F values
# for rows whose current values match the "before" snapshot, in both
# the "展览" sheet and the "全部类型" sheet.

$wb = $excel.ActiveWorkbook

# Map of old value -> new value for column F, applied per worksheet/row.
$updates = @(
    @{ Sheet = "展览"; Row = 3;  Old = 5623; New = 5626 },
    @{ Sheet = "展览"; Row = 4;  Old = 79;   New = 80 },
    @{ Sheet = "展览"; Row = 6;  Old = 936;  New = 938 },
    @{ Sheet = "展览"; Row = 8;  Old = 2517; New = 2524 },
    @{ Sheet = "展览"; Row = 10; Old = 146;  New = 148 },
    @{ Sheet = "展览"; Row = 11; Old = 9;    New = 10 },
    @{ Sheet = "展览"; Row = 13; Old = 18;   New = 19 },
    @{ Sheet = "展览"; Row = 14; Old = 2362; New = 2363 },
    @{ Sheet = "展览"; Row = 15; Old = 350;  New = 356 },

    @{ Sheet = "全部类型"; Row = 3;  Old = 5623; New = 5626 },
    @{ Sheet = "全部类型"; Row = 5;  Old = 79;   New = 80 },
    @{ Sheet = "全部类型"; Row = 8;  Old = 936;  New = 938 },
    @{ Sheet = "全部类型"; Row = 10; Old = 2517; New = 2524 },
    @{ Sheet = "全部类型"; Row = 12; Old = 146;  New = 148 },
    @{ Sheet = "全部类型"; Row = 13; Old = 9;    New = 10 },
    @{ Sheet = "全部类型"; Row = 16; Old = 18;   New = 19 },
    @{ Sheet = "全部类型"; Row = 17; Old = 2362; New = 2363 },
    @{ Sheet = "全部类型"; Row = 18; Old = 350;  New = 356 }
)

foreach ($u in $updates) {
    $ws = $wb.Worksheets.Item($u.Sheet)
    $cell = $ws.Cells.Item($u.Row, 6)  # Column F is the 6th column
    $current = $cell.Value()
    if ($current -ne $u.Old) {
        Write-Host "Warning: $($u.Sheet)!F$($u.Row) expected $($u.Old) but found $current"
    }
    $cell.Value = $u.New
}

Write-Host "Updated $($updates.Count) cells across sheets 展览 and 全部类型."
